$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1409
$ws.Cells.Item(3, 6).Value = 102
$ws.Cells.Item(4, 6).Value = 2108
$ws.Cells.Item(5, 6).Value = 6297
$ws.Cells.Item(6, 6).Value = 500
$ws.Cells.Item(7, 6).Value = 1052
$ws.Cells.Item(9, 6).Value = 3471
$ws.Cells.Item(10, 6).Value = 6704
$ws.Cells.Item(11, 6).Value = 211
$ws.Cells.Item(12, 6).Value = 1343
$ws.Cells.Item(13, 6).Value = 779
$ws.Cells.Item(15, 6).Value = 10
$ws.Cells.Item(16, 6).Value = 27
$ws.Cells.Item(17, 6).Value = 1123
$ws.Cells.Item(19, 6).Value = 114
$ws.Cells.Item(23, 6).Value = 1002
$ws.Cells.Item(24, 6).Value = 325
$ws.Cells.Item(25, 6).Value = 34
$ws.Cells.Item(26, 6).Value = 20
$ws.Cells.Item(27, 6).Value = 111
$ws.Cells.Item(30, 6).Value = 22
$ws.Cells.Item(31, 6).Value = 69
$ws.Cells.Item(33, 6).Value = 22
$ws.Cells.Item(34, 6).Value = 22
$ws.Cells.Item(35, 6).Value = 5
$ws.Cells.Item(36, 6).Value = 321
$ws.Cells.Item(39, 6).Value = 297
$ws.Cells.Item(40, 6).Value = 1176

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(11, 6).Value = 116
$ws.Cells.Item(12, 6).Value = 13
$ws.Cells.Item(16, 6).Value = 1722
$ws.Cells.Item(17, 6).Value = 547
$ws.Cells.Item(19, 6).Value = 7
$ws.Cells.Item(21, 6).Value = 185
$ws.Cells.Item(30, 6).Value = 695
$ws.Cells.Item(31, 6).Value = 955
$ws.Cells.Item(32, 6).Value = 573
$ws.Cells.Item(34, 6).Value = 86

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 703
$ws.Cells.Item(5, 6).Value = 829
$ws.Cells.Item(6, 6).Value = 574
$ws.Cells.Item(8, 6).Value = 1124

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 703
$ws.Cells.Item(3, 6).Value = 829
$ws.Cells.Item(5, 6).Value = 102
$ws.Cells.Item(6, 6).Value = 574
$ws.Cells.Item(10, 6).Value = 6297
$ws.Cells.Item(11, 6).Value = 500
$ws.Cells.Item(12, 6).Value = 1052
$ws.Cells.Item(15, 6).Value = 6704
$ws.Cells.Item(16, 6).Value = 116
$ws.Cells.Item(17, 6).Value = 211
$ws.Cells.Item(18, 6).Value = 1343
$ws.Cells.Item(19, 6).Value = 13
$ws.Cells.Item(22, 6).Value = 1722
$ws.Cells.Item(23, 6).Value = 547
$ws.Cells.Item(24, 6).Value = 1124
$ws.Cells.Item(25, 6).Value = 10
$ws.Cells.Item(26, 6).Value = 7
$ws.Cells.Item(27, 6).Value = 185
$ws.Cells.Item(29, 6).Value = 114
$ws.Cells.Item(32, 6).Value = 34
$ws.Cells.Item(33, 6).Value = 20
$ws.Cells.Item(36, 6).Value = 22
$ws.Cells.Item(37, 6).Value = 69
$ws.Cells.Item(40, 6).Value = 955
$ws.Cells.Item(41, 6).Value = 22
$ws.Cells.Item(42, 6).Value = 573
$ws.Cells.Item(43, 6).Value = 321
$ws.Cells.Item(45, 6).Value = 86
$ws.Cells.Item(46, 6).Value = 297
